$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "{'label': 'LABEL_0', 'score': 0.6165294051170349}"
$ws.Range("H3").Value = "{'label': 'LABEL_0', 'score': 0.6152240633964539}"
$ws.Range("H4").Value = "{'label': 'LABEL_0', 'score': 0.6276030540466309}"
$ws.Range("H5").Value = "{'label': 'LABEL_0', 'score': 0.6287968158721924}"
$ws.Range("H6").Value = "{'label': 'LABEL_0', 'score': 0.6265703439712524}"
$ws.Range("H7").Value = "{'label': 'LABEL_0', 'score': 0.6092011332511902}"
$ws.Range("H8").Value = "{'label': 'LABEL_0', 'score': 0.6154984831809998}"
$ws.Range("H9").Value = "{'label': 'LABEL_0', 'score': 0.6121871471405029}"
$ws.Range("H10").Value = "{'label': 'LABEL_0', 'score': 0.6297902464866638}"
$ws.Range("H11").Value = "{'label': 'LABEL_0', 'score': 0.6118900179862976}"
$ws.Range("H12").Value = "{'label': 'LABEL_0', 'score': 0.6316757798194885}"
$ws.Range("H13").Value = "{'label': 'LABEL_0', 'score': 0.6154239177703857}"
$ws.Range("H14").Value = "{'label': 'LABEL_0', 'score': 0.6156269311904907}"
$ws.Range("H15").Value = "{'label': 'LABEL_0', 'score': 0.623571515083313}"
$ws.Range("H16").Value = "{'label': 'LABEL_0', 'score': 0.6084762215614319}"
$ws.Range("H17").Value = "{'label': 'LABEL_0', 'score': 0.5984631776809692}"
$ws.Range("H18").Value = "{'label': 'LABEL_0', 'score': 0.6161437034606934}"
$ws.Range("H19").Value = "{'label': 'LABEL_0', 'score': 0.6115897297859192}"
$ws.Range("H20").Value = "{'label': 'LABEL_0', 'score': 0.6084784269332886}"
$ws.Range("H21").Value = "{'label': 'LABEL_0', 'score': 0.6260745525360107}"
$ws.Range("H22").Value = "{'label': 'LABEL_0', 'score': 0.6434668302536011}"
$ws.Range("H23").Value = "{'label': 'LABEL_0', 'score': 0.626284658908844}"
$ws.Range("H24").Value = "{'label': 'LABEL_0', 'score': 0.612806499004364}"
$ws.Range("H25").Value = "{'label': 'LABEL_0', 'score': 0.6143962740898132}"
$ws.Range("H26").Value = "{'label': 'LABEL_0', 'score': 0.6199201345443726}"
$ws.Range("H27").Value = "{'label': 'LABEL_0', 'score': 0.6199201345443726}"
$ws.Range("H28").Value = "{'label': 'LABEL_0', 'score': 0.6308311820030212}"
$ws.Range("H29").Value = "{'label': 'LABEL_0', 'score': 0.6215488910675049}"
$ws.Range("H30").Value = "{'label': 'LABEL_0', 'score': 0.644977867603302}"
$ws.Range("H31").Value = "{'label': 'LABEL_0', 'score': 0.6224562525749207}"
$ws.Range("H32").Value = "{'label': 'LABEL_0', 'score': 0.6275407075881958}"
$ws.Range("H33").Value = "{'label': 'LABEL_0', 'score': 0.6275407075881958}"
$ws.Range("H34").Value = "{'label': 'LABEL_0', 'score': 0.6275407075881958}"
$ws.Range("H35").Value = "{'label': 'LABEL_0', 'score': 0.6146340370178223}"
$ws.Range("H36").Value = "{'label': 'LABEL_0', 'score': 0.6022341847419739}"
$ws.Range("H37").Value = "{'label': 'LABEL_0', 'score': 0.6090781092643738}"
$ws.Range("H38").Value = "{'label': 'LABEL_0', 'score': 0.6110119819641113}"
$ws.Range("H39").Value = "{'label': 'LABEL_0', 'score': 0.648396909236908}"
$ws.Range("H40").Value = "{'label': 'LABEL_0', 'score': 0.6297367811203003}"
$ws.Range("H41").Value = "{'label': 'LABEL_0', 'score': 0.62285977602005}"
$ws.Range("H42").Value = "{'label': 'LABEL_0', 'score': 0.6285818219184875}"
$ws.Range("H43").Value = "{'label': 'LABEL_0', 'score': 0.6175370216369629}"
$ws.Range("H44").Value = "{'label': 'LABEL_0', 'score': 0.6186134815216064}"
$ws.Range("H45").Value = "{'label': 'LABEL_0', 'score': 0.6211050152778625}"
$ws.Range("H46").Value = "{'label': 'LABEL_0', 'score': 0.612415611743927}"
$ws.Range("H47").Value = "{'label': 'LABEL_0', 'score': 0.6200401782989502}"
$ws.Range("H48").Value = "{'label': 'LABEL_0', 'score': 0.6257834434509277}"
$ws.Range("H49").Value = "{'label': 'LABEL_0', 'score': 0.6227863430976868}"
$ws.Range("H50").Value = "{'label': 'LABEL_0', 'score': 0.6207936406135559}"
$ws.Range("H51").Value = "{'label': 'LABEL_0', 'score': 0.6283147931098938}"
$ws.Range("H52").Value = "{'label': 'LABEL_0', 'score': 0.6177467703819275}"
$ws.Range("H53").Value = "{'label': 'LABEL_0', 'score': 0.618755578994751}"
$ws.Range("H54").Value = "{'label': 'LABEL_0', 'score': 0.6097444295883179}"
$ws.Range("H55").Value = "{'label': 'LABEL_0', 'score': 0.6146154999732971}"
$ws.Range("H56").Value = "{'label': 'LABEL_0', 'score': 0.6213558912277222}"
$ws.Range("H57").Value = "{'label': 'LABEL_0', 'score': 0.6223445534706116}"
$ws.Range("H58").Value = "{'label': 'LABEL_0', 'score': 0.6208321452140808}"
$ws.Range("H59").Value = "{'label': 'LABEL_0', 'score': 0.649174153804779}"
$ws.Range("H60").Value = "{'label': 'LABEL_0', 'score': 0.59809410572052}"
$ws.Range("H61").Value = "{'label': 'LABEL_0', 'score': 0.6069620847702026}"
$ws.Range("H62").Value = "{'label': 'LABEL_0', 'score': 0.6208321452140808}"
$ws.Range("H63").Value = "{'label': 'LABEL_0', 'score': 0.6127477884292603}"
$ws.Range("H64").Value = "{'label': 'LABEL_0', 'score': 0.6205687522888184}"
$ws.Range("H65").Value = "{'label': 'LABEL_0', 'score': 0.6281725168228149}"
$ws.Range("H66").Value = "{'label': 'LABEL_0', 'score': 0.6178492307662964}"
$ws.Range("H67").Value = "{'label': 'LABEL_0', 'score': 0.6152944564819336}"
$ws.Range("H68").Value = "{'label': 'LABEL_0', 'score': 0.6149681210517883}"
$ws.Range("H69").Value = "{'label': 'LABEL_0', 'score': 0.622868001461029}"
$ws.Range("H70").Value = "{'label': 'LABEL_0', 'score': 0.6356972455978394}"
$ws.Range("H71").Value = "{'label': 'LABEL_0', 'score': 0.6169376373291016}"
$ws.Range("H72").Value = "{'label': 'LABEL_0', 'score': 0.6011926531791687}"
$ws.Range("H73").Value = "{'label': 'LABEL_0', 'score': 0.6514390707015991}"
$ws.Range("H74").Value = "{'label': 'LABEL_0', 'score': 0.6188256740570068}"
$ws.Range("H75").Value = "{'label': 'LABEL_0', 'score': 0.6082078218460083}"
$ws.Range("H76").Value = "{'label': 'LABEL_0', 'score': 0.6253708600997925}"
$ws.Range("H77").Value = "{'label': 'LABEL_0', 'score': 0.6273909211158752}"
$ws.Range("H78").Value = "{'label': 'LABEL_0', 'score': 0.6172049045562744}"
$ws.Range("H79").Value = "{'label': 'LABEL_0', 'score': 0.6198626756668091}"
$ws.Range("H80").Value = "{'label': 'LABEL_0', 'score': 0.6217631101608276}"
$ws.Range("H81").Value = "{'label': 'LABEL_0', 'score': 0.6265097856521606}"
$ws.Range("H82").Value = "{'label': 'LABEL_0', 'score': 0.6118432879447937}"
$ws.Range("H83").Value = "{'label': 'LABEL_0', 'score': 0.6141948103904724}"
$ws.Range("H84").Value = "{'label': 'LABEL_0', 'score': 0.6261197924613953}"
$ws.Range("H85").Value = "{'label': 'LABEL_0', 'score': 0.6199622750282288}"
$ws.Range("H86").Value = "{'label': 'LABEL_0', 'score': 0.6402068138122559}"
$ws.Range("H87").Value = "{'label': 'LABEL_0', 'score': 0.612963080406189}"
$ws.Range("H88").Value = "{'label': 'LABEL_0', 'score': 0.6144188046455383}"
$ws.Range("H89").Value = "{'label': 'LABEL_0', 'score': 0.613923966884613}"
$ws.Range("H90").Value = "{'label': 'LABEL_0', 'score': 0.6049232482910156}"
$ws.Range("H91").Value = "{'label': 'LABEL_0', 'score': 0.6211430430412292}"
$ws.Range("H92").Value = "{'label': 'LABEL_0', 'score': 0.6153357028961182}"
$ws.Range("H93").Value = "{'label': 'LABEL_0', 'score': 0.5993615388870239}"
$ws.Range("H94").Value = "{'label': 'LABEL_0', 'score': 0.6199654936790466}"
$ws.Range("H95").Value = "{'label': 'LABEL_0', 'score': 0.6265184879302979}"
$ws.Range("H96").Value = "{'label': 'LABEL_0', 'score': 0.6110716462135315}"
$ws.Range("H97").Value = "{'label': 'LABEL_0', 'score': 0.6054821610450745}"
$ws.Range("H98").Value = "{'label': 'LABEL_0', 'score': 0.6227136254310608}"
$ws.Range("H99").Value = "{'label': 'LABEL_0', 'score': 0.6235061287879944}"
$ws.Range("H100").Value = "{'label': 'LABEL_0', 'score': 0.6221509575843811}"
$ws.Range("H101").Value = "{'label': 'LABEL_0', 'score': 0.6235061287879944}"
